$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths ---
# NOTE: the ColumnWidth object-model property (in characters) gets a fixed
# +5/6 padding baked into the stored OOXML "width" attribute by this engine
# (mirrors real Excel's MDW-based width storage). Subtract 5/6 from the
# desired stored width so the resulting file matches exactly.
$ws.Columns.Item(3).ColumnWidth = 52 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 44 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 15 - (5/6)
$ws.Columns.Item(7).ColumnWidth = 15 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 34 - (5/6)

# --- New data rows 2-7 ---
$data = @(
    @("1330036", "https://aiesec.org/opportunity/global-talent/1330036", "[Accelerate Serbia] Structural engineer", "Belgrade, Serbia", "No", "0 applicants", "9 - 12 Weeks", "Welt Inzenjering"),
    @("1330033", "https://aiesec.org/opportunity/global-talent/1330033", "[Accelerate Serbia] Design development assistance", "Belgrade, Serbia", "No", "0 applicants", "9 - 12 Weeks", "IWA CONSALT DOO"),
    @("1329992", "https://aiesec.org/opportunity/global-talent/1329992", "Business Development manager", "Nugegoda, Sri Lanka", "No", "0 applicants", "3 - 6 Months", "Weblook International (Pvt) Ltd"),
    @("1329871", "https://aiesec.org/opportunity/global-talent/1329871", "Data Analyst Intern", "Mayur Vihar, Delhi, India", "No", "0 applicants", "3 - 6 Months", "Credifin Limited"),
    @("1329869", "https://aiesec.org/opportunity/global-talent/1329869", "Marketing and Strategy Intern", "Mayur Vihar, Delhi, India", "No", "0 applicants", "3 - 6 Months", "Credifin Limited"),
    @("1327809", "https://aiesec.org/opportunity/global-talent/1327809", "Video editor", "El Sadat City, Menofia Governorate, Egypt", "No", "2 applicants", "9 - 12 Weeks", "Habib Agency")
)

$rowIndex = 2
foreach ($rowData in $data) {
    for ($col = 1; $col -le 8; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col - 1]
        $cell.ClearFormats()
    }
    $rowIndex++
}
